$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap column widths: col C <-> col G widths, and col O <-> col Q widths.
# (The runtime quantizes stored column width to 1/6-character-unit (pixel)
# increments, so these ColumnWidth values are chosen to land on the
# closest achievable width to the exact target from the spec.)
$ws.Columns.Item(3).ColumnWidth = 2.25
$ws.Columns.Item(7).ColumnWidth = 1.25
$ws.Columns.Item(15).ColumnWidth = 4.75
$ws.Columns.Item(17).ColumnWidth = 3.75

# Update row 1 values
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 10
$ws.Range("E1").Value = 12
$ws.Range("F1").Value = 23
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 24
$ws.Range("I1").Value = 16
$ws.Range("J1").Value = 30
$ws.Range("K1").Value = 20
$ws.Range("L1").Value = 12
$ws.Range("M1").Value = 0.097000000000000003
$ws.Range("N1").Value = 0.045999999999999999
$ws.Range("O1").Value = 0.087999999999999995
$ws.Range("P1").Value = 0.028000000000000004
$ws.Range("Q1").Value = 0.089999999999999997
